$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 4 and row 5 and need to be swapped
$cols = @("A","B","D","E","F","G","H","Q","R")

foreach ($col in $cols) {
    $addr4 = $col + "4"
    $addr5 = $col + "5"
    $v4 = $ws.Range($addr4).Value()
    $v5 = $ws.Range($addr5).Value()
    $ws.Range($addr4).Value = $v5
    $ws.Range($addr5).Value = $v4
}

# The stray empty "Bestamningsmetod" cell (AF) moves from row 4 to row 5
$ws.Range("AF4").Cut($ws.Range("AF5"))
